$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-height tweaks: 18.75 -> 19.5 for the header row and the first few data rows
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(5).RowHeight = 19.5
$ws.Rows.Item(6).RowHeight = 19.5
$ws.Rows.Item(7).RowHeight = 19.5

# "Binary Search" Python attempt count bumped from 3 to 4
$ws.Range("E10").Value = 4

# New entry: "Delete Node in a Linked List"
$ws.Range("A33").Value = "LC"
$ws.Range("B33").Value = "Delete Node in a Linked List"
$ws.Range("C33").Value = "Medium"
$ws.Range("D33").Value = "NA"

# E33 needs the same number format/alignment/font as the other attempt-count
# cells (style index 8) rather than the blank-row default (style index 3) -
# copy formats from an already-styled attempt-count cell, then set the value.
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = 1

$ws.Range("K33").Value = "O(1)"

$excel.CutCopyMode = 0
